$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): translate column headers to snake_case names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize capitalization of particles ("de" -> "De", "del" -> "Del", "los" -> "Los", "y" -> "Y") ---
$ws.Range("B8").Value  = "Amatenango De La Frontera"
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("A27").Value = "Estado De México"
$ws.Range("B30").Value = "San Felipe Del Progreso"
$ws.Range("B37").Value = "Acapulco De Juárez"
$ws.Range("B39").Value = "Alcozauca De Guerrero"
$ws.Range("B40").Value = "Atenango Del Río"
$ws.Range("B42").Value = "Atoyac De Álvarez"
$ws.Range("B43").Value = "Chilapa De Álvarez"
$ws.Range("B44").Value = "Chilpancingo De Los Bravo"
$ws.Range("B50").Value = "Tlapa De Comonfort"
$ws.Range("A59").Value = "Michoacán De Ocampo"
$ws.Range("B74").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B79").Value = "Los Reyes De Juárez"
$ws.Range("B83").Value = "Tepanco De López"
$ws.Range("B85").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B86").Value = "Xayacatlán De Bravo"
$ws.Range("A107").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B109").Value = "Amatlán De Los Reyes"
$ws.Range("A115").Value = "Total"

# --- Remove trailing footer/source note rows (117-121) ---
$ws.Range("A117:A121").EntireRow.Delete()
